$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header in column F, row 4
$ws.Range("F4").Value = "Tanggal Jatuh Tempo"

# Set column F width to match the diff (target stored width 20.6328125 chars).
# NOTE: the host's ColumnWidth setter quantizes to a 1/6-character pixel grid
# (xml_width = (round(ColumnWidth*6)+5)/6), so 19.83 is the closest input that
# lands on the nearest achievable stored width (20.666666...) to the target.
$ws.Columns.Item(6).ColumnWidth = 19.83

# Update the selected cell to F10 (matches the diff's <selection activeCell="F10" sqref="F10"/>)
$ws.Range("F10").Select()
